$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 325 and 326 (columns C-F)
$ws.Range("C325").Value = 7710432000000
$ws.Range("D325").Value = 7710432000000
$ws.Range("E325").Value = 7710432000000
$ws.Range("F325").Value = 7710432000000

$ws.Range("C326").Value = 7822810000000
$ws.Range("D326").Value = 7822810000000
$ws.Range("E326").Value = 7822810000000
$ws.Range("F326").Value = 7822810000000

# Add new rows 327-329
$ws.Range("A327").Value = 44986.45833333334
$ws.Range("B327").Value = "ECONOMICS:EGM2"
$ws.Range("C327").Value = 7965088000000
$ws.Range("D327").Value = 7965088000000
$ws.Range("E327").Value = 7965088000000
$ws.Range("F327").Value = 7965088000000
$ws.Range("G327").Value = 0

$ws.Range("A328").Value = 45017.45833333334
$ws.Range("B328").Value = "ECONOMICS:EGM2"
$ws.Range("C328").Value = 8069151000000
$ws.Range("D328").Value = 8069151000000
$ws.Range("E328").Value = 8069151000000
$ws.Range("F328").Value = 8069151000000
$ws.Range("G328").Value = 0

$ws.Range("A329").Value = 45047.41666666666
$ws.Range("B329").Value = "ECONOMICS:EGM2"
$ws.Range("C329").Value = 8140535000000
$ws.Range("D329").Value = 8140535000000
$ws.Range("E329").Value = 8140535000000
$ws.Range("F329").Value = 8140535000000
$ws.Range("G329").Value = 0

# Apply same style as other date cells (A2 style) to new rows
$ws.Range("A2").Copy()
$ws.Range("A327:A329").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
